# edit.ps1
# Applies the "New crime data collected" update to the CompStat weekly report.
#
# Summary of changes:
#  - Report header: Volume Number 35 -> 36
#  - Report header: week covering 8/28/2023 - 9/3/2023 -> 9/4/2023 - 9/10/2023
#  - Crime Complaints table (rows 15-30): refreshed counts and computed % changes

function Set-TextZero($ws, $cellRef) {
    # Some cells in this report represent a value of 0 as the literal text "0"
    # (general format, right aligned) instead of a formatted number.
    # A leading apostrophe forces Excel to store the value as text.
    $ws.Range($cellRef).Value = "'0"
}

function Set-NumericFromText($ws, $cellRef, $num) {
    # Converts a cell that currently holds the text "0" back into a real,
    # "#,##0"-formatted number.
    $ws.Range($cellRef).NumberFormat = "#,##0"
    $ws.Range($cellRef).Value = $num
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (rich-text shared strings)
# ---------------------------------------------------------------------------

# "Volume 30   Number  35" -> "Volume 30   Number  36"
$ws.Range("A8").Characters(21,2).Text = "36"

# "Report Covering the Week  8/28/2023  Through  9/3/2023"
#   -> "Report Covering the Week  9/4/2023  Through  9/10/2023"
$ws.Range("C9").Characters(27,9).Text = "9/4/2023"
$ws.Range("C9").Characters(46,8).Text = "9/10/2023"

# ---------------------------------------------------------------------------
# Crime Complaints table updates
# ---------------------------------------------------------------------------

Set-TextZero $ws "F15"
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -100
$ws.Range("J15").Value = 17
$ws.Range("K15").Value = -41.176470588235
$ws.Range("M15").Value = 11.111111111111
$ws.Range("N15").Value = -68.75
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -37.5
$ws.Range("F16").Value = 34
$ws.Range("G16").Value = 33
$ws.Range("H16").Value = 3.030303030303
$ws.Range("I16").Value = 332
$ws.Range("J16").Value = 420
$ws.Range("K16").Value = -20.952380952381
$ws.Range("L16").Value = 24.344569288389
$ws.Range("M16").Value = 201.818181818182
$ws.Range("N16").Value = -81.586245146977
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = -18.181818181818
$ws.Range("F17").Value = 47
$ws.Range("G17").Value = 41
$ws.Range("H17").Value = 14.634146341463
$ws.Range("I17").Value = 352
$ws.Range("J17").Value = 310
$ws.Range("K17").Value = 13.548387096774
$ws.Range("L17").Value = 7.975460122699
$ws.Range("M17").Value = 158.823529411765
$ws.Range("N17").Value = -24.786324786324
$ws.Range("C18").Value = 12
$ws.Range("D18").Value = 13
$ws.Range("E18").Value = -7.692307692307
$ws.Range("F18").Value = 25
$ws.Range("G18").Value = 52
$ws.Range("H18").Value = -51.923076923076
$ws.Range("I18").Value = 292
$ws.Range("J18").Value = 474
$ws.Range("K18").Value = -38.396624472573
$ws.Range("L18").Value = 2.456140350877
$ws.Range("M18").Value = 28.63436123348
$ws.Range("N18").Value = -84.284176533907
$ws.Range("C19").Value = 36
$ws.Range("D19").Value = 48
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 156
$ws.Range("G19").Value = 188
$ws.Range("H19").Value = -17.021276595744
$ws.Range("I19").Value = 1610
$ws.Range("J19").Value = 1547
$ws.Range("K19").Value = 4.072398190045
$ws.Range("L19").Value = 82.746878547105
$ws.Range("M19").Value = 2.744097000638
$ws.Range("N19").Value = -75.6687320538
$ws.Range("C20").Value = 3
$ws.Range("F20").Value = 9
$ws.Range("H20").Value = 800
$ws.Range("I20").Value = 49
$ws.Range("K20").Value = 2.083333333333
$ws.Range("L20").Value = 44.117647058823
$ws.Range("M20").Value = 226.666666666667
$ws.Range("N20").Value = -80.321285140562
$ws.Range("C21").Value = 65
$ws.Range("D21").Value = 81
$ws.Range("E21").Value = -19.753086419753
$ws.Range("F21").Value = 271
$ws.Range("G21").Value = 319
$ws.Range("H21").Value = -15.047021943573
$ws.Range("I21").Value = 2647
$ws.Range("J21").Value = 2822
$ws.Range("K21").Value = -6.201275690999
$ws.Range("L21").Value = 46.485888212506
$ws.Range("M21").Value = 28.246124031007
$ws.Range("N21").Value = -76.014860456687
Set-NumericFromText $ws "C22" 2
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 0
$ws.Range("G22").Value = 12
$ws.Range("H22").Value = -8.333333333333
$ws.Range("I22").Value = 144
$ws.Range("J22").Value = 126
$ws.Range("K22").Value = 14.285714285714
$ws.Range("L22").Value = 42.574257425742
$ws.Range("M22").Value = 46.938775510204
$ws.Range("C24").Value = 67
$ws.Range("D24").Value = 68
$ws.Range("E24").Value = -1.470588235294
$ws.Range("F24").Value = 316
$ws.Range("G24").Value = 302
$ws.Range("H24").Value = 4.635761589403
$ws.Range("I24").Value = 2816
$ws.Range("J24").Value = 2270
$ws.Range("K24").Value = 24.052863436123
$ws.Range("L24").Value = 93.406593406593
$ws.Range("M24").Value = -13.327177593105
$ws.Range("C25").Value = 20
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 90
$ws.Range("G25").Value = 60
$ws.Range("H25").Value = 50
$ws.Range("I25").Value = 752
$ws.Range("J25").Value = 612
$ws.Range("K25").Value = 22.875816993464
$ws.Range("L25").Value = 23.076923076923
$ws.Range("M25").Value = 87.531172069825
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = -80
$ws.Range("J26").Value = 22
$ws.Range("K26").Value = -13.636363636363
Set-TextZero $ws "C27"
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 15
$ws.Range("G27").Value = 22
$ws.Range("H27").Value = -31.818181818181
$ws.Range("J27").Value = 158
$ws.Range("K27").Value = -1.898734177215
$ws.Range("L27").Value = 47.619047619047
$ws.Range("L30").Value = -65.384615384615
